$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix letter case inconsistency: "Vocabulary Code" -> "Vocabulary code"
$ws.Range("H2").Value = "Vocabulary code"

# Reflect the active cell selection left after the edit
$ws.Range("H2").Select()
